$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert two rows before row 11. This pushes the existing "columnsToKeep"
# setting (currently row 11) down to row 13, leaving row 12 blank exactly
# like the untouched row 12 already was, and it naturally shifts every
# following populated/blank row down by two (so the sheet grows from 994
# to 996 rows, matching the new dimension).
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# Populate the brand-new row 11 with the "companiiFilePath" setting.
# (Value/Description are entered before the Name, matching the order the
# new shared strings were appended in the authored workbook.)
$ws.Cells.Item(11,2).Value = "Data\Output\companii.xlsx"
$ws.Cells.Item(11,3).Value = 'Location to save the file "companii"'
$ws.Cells.Item(11,1).Value = "companiiFilePath"

# Row 11 uses the plain (non-wrapped) row height/format used by the other
# simple Name/Value/Description rows such as row 9.
$ws.Rows.Item(11).RowHeight = 14.25
$ws.Rows.Item(12).RowHeight = 14.25

# Update the active selection to the single cell A11 (was the whole row).
$ws.Range("A11").Select()
